$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the legacy cell-comments (A3..T3, except G3 which has none).
#    Deleting every Comment object removes comments1.xml / vmlDrawing1.vml
#    and the legacyDrawing reference on the sheet, same as the target.
# ---------------------------------------------------------------------------
$commentCols = @("A","B","C","D","E","F","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($col in $commentCols) {
    $cell = $ws.Range($col + "3")
    $cm = $cell.Comment
    if ($cm -ne $null) {
        $cm.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Shift the header row up from row 3 to row 2 by deleting the blank
#    row 2 above it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 3. Populate the new row 3 with the descriptive text that used to live in
#    the cell comments (plus a couple of brand-new notes).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "eg. 1999"
$ws.Range("B3").Value = "eg, Apr, Sep"
$ws.Range("C3").Value = "eg. 1"
$ws.Range("D3").Value = "Must match site name in database"
$ws.Range("E3").Value = "Eg. Bonell. OptionalOptional Y/N value"
$ws.Range("F3").Value = "Fill only if fish has no pit tag. Value can be arbritrary (eg. 1, 2, 3), but must be unique per fish."
$ws.Range("G3").Value = "Pit tag number of fist."
$ws.Range("H3").Value = "Optional. Was a new pit tag used? Y/N"
$ws.Range("I3").Value = "Optional. Must match animal subjective detail in database. Eg. Bonell. "
$ws.Range("J3").Value = "Eg. FP 2021"
$ws.Range("K3").Value = "Optional. Used to indicate if fish were brought back to the facility or returned to the river."
$ws.Range("L3").Value = "Optional. Units can be set in header to (cm) or (mm)."
$ws.Range("M3").Value = "Optional. Units can be set in header to (g) or (kg)."
$ws.Range("N3").Value = "eg. M/F/I"
$ws.Range("O3").Value = "Optional. Y/N"
$ws.Range("P3").Value = "Vial Number, optional"
$ws.Range("Q3").Value = "Optional. Y/N."
$ws.Range("R3").Value = "Optional. Y/N"
$ws.Range("S3").Value = "Initials of team at site. Eg. AB, CD"
$ws.Range("T3").Value = "Optional"

# ---------------------------------------------------------------------------
# 4. Style row 3: small grey Tahoma font, top-aligned + wrapped, boxed in
#    thin borders. Build the format once on a scratch cell and paste it in
#    a single shot so we don't litter the style table with partial states.
# ---------------------------------------------------------------------------
$helper = $ws.Range("Z1")
$helper.Font.Name = "Tahoma"
$helper.Font.Size = 10
$helper.Font.Color = 8421504
$helper.WrapText = $true
$helper.VerticalAlignment = -4160
$helper.Borders.LineStyle = 1

$helper.Copy()
$target = $ws.Range("A3:T3")
$target.PasteSpecial(-4122, $null, $null, $null)
$helper.Clear()
$excel.CutCopyMode = 0

$ws.Rows.Item(3).RowHeight = 63.75

# ---------------------------------------------------------------------------
# 5. Column width tweaks so the longer note text is legible.
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 22.2573
$ws.Columns.Item(9).ColumnWidth = 16.687
$ws.Columns.Item(11).ColumnWidth = 20.1167
$ws.Columns.Item(12).ColumnWidth = 12.9721
$ws.Columns.Item(13).ColumnWidth = 13.8276
$ws.Columns.Item(14).ColumnWidth = 20.8276
$ws.Columns.Item(16).ColumnWidth = 11.687
$ws.Columns.Item(17).ColumnWidth = 10.4168

# ---------------------------------------------------------------------------
# 6. Selection cosmetics, matching the saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("A1:T2").Select()
